# sesDesignDoc.pptx edit: "Character work, story work, greyboxing start"
#
# 1) Retitle the working-title textbox on slide 1 from
#      "Someone Else's Shoes (Working Title)"
#    to
#      "I Was Here (Working Title)"
#    (split across two runs, matching how PowerPoint re-runs text after
#    an in-place partial edit).
# 2) Append a new, blank "Title and Content" slide (slide 2) to the deck.

$p = $ppt.ActivePresentation

# --- 1) Update the title textbox on slide 1 -------------------------------
$s1 = $p.Slides.Item(1)
$titleBox = $s1.Shapes.Item(2)   # "TextBox 4" - the working-title caption

$titleBox.TextFrame.TextRange.Text = "I Was Here (Working Title)"
# Re-touch the leading portion so it becomes its own run (mirrors the
# authored diff, which shows the text split into two runs with identical
# formatting).
$lead = $titleBox.TextFrame.TextRange.Characters(1, 12)
$lead.Text = "I Was Here ("

# --- 2) Add a new slide (Title and Content layout) ------------------------
$newSlide = $p.Slides.Add(2, 2)
